$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.323.76"
$ws.Range("E2").Value = "  +2.43%  "

$ws.Range("D3").Value = "3.204.07"
$ws.Range("E3").Value = "  +1.93%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'597.99"
$ws.Range("E5").Value = "  +1.80%  "

$ws.Range("D6").Value = "'154.09"
$ws.Range("E6").Value = "  +5.82%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "3.201.32"
$ws.Range("E8").Value = "  +1.97%  "

$ws.Range("D9").Value = "'0.543"
$ws.Range("E9").Value = "  +2.51%  "

$ws.Range("D10").Value = "'0.168"
$ws.Range("E10").Value = "  +4.25%  "

$ws.Range("D11").Value = "'6.10"
$ws.Range("E11").Value = "  +6.10%  "

$ws.Range("E12").Value = "  +2.81%  "

$ws.Range("E13").Value = "  +3.36%  "

$ws.Range("D14").Value = "'39.47"
$ws.Range("E14").Value = "  +6.55%  "

$ws.Range("D15").Value = "3.729.79"
$ws.Range("E15").Value = "  +1.80%  "

$ws.Range("E16").Value = "  +0.40%  "

$ws.Range("D17").Value = "'7.43"
$ws.Range("E17").Value = "  +4.63%  "

$ws.Range("D18").Value = "65.044.09"
$ws.Range("E18").Value = "  +2.30%  "

$ws.Range("D19").Value = "3.201.07"
$ws.Range("E19").Value = "  +1.79%  "

$ws.Range("D20").Value = "'483.33"
$ws.Range("E20").Value = "  +4.28%  "

$ws.Range("D21").Value = "'15.15"
$ws.Range("E21").Value = "  +6.15%  "

$ws.Range("D22").Value = "'0.774"
$ws.Range("E22").Value = "  +6.02%  "

$ws.Range("D23").Value = "'7.93"
$ws.Range("E23").Value = "  +6.57%  "

$ws.Range("D24").Value = "'13.97"
$ws.Range("E24").Value = "  +7.55%  "

$ws.Range("D25").Value = "'2.44"
$ws.Range("E25").Value = "  +11.01%  "

$ws.Range("D26").Value = "'83.69"
$ws.Range("E26").Value = "  +2.88%  "

$ws.Range("E27").Value = "  +0.31%  "

$ws.Range("D28").Value = "'9.90"
$ws.Range("E28").Value = "  +8.64%  "

$ws.Range("E29").Value = "  +3.70%  "

$ws.Range("D30").Value = "'7.52"
$ws.Range("E30").Value = "  +7.98%  "

$ws.Range("E31").Value = "  +3.17%  "

$ws.Range("E32").Value = "  +0.06%  "

$ws.Range("E33").Value = "  +9.17%  "

$ws.Range("D34").Value = "'28.58"
$ws.Range("E34").Value = "  +5.83%  "

$ws.Range("D35").Value = "0.0₃0902"
$ws.Range("E35").Value = "  +5.60%  "

$ws.Range("D36").Value = "'3.59"
$ws.Range("E36").Value = "  +6.75%  "

$ws.Range("E37").Value = "  +4.63%  "

$ws.Range("E38").Value = "  +5.79%  "

$ws.Range("E39").Value = "  +3.50%  "

$ws.Range("D40").Value = "'479.32"
$ws.Range("E40").Value = "  +8.95%  "

$ws.Range("E41").Value = "  +7.95%  "

$ws.Range("D42").Value = "'51.45"
$ws.Range("E42").Value = "  +1.18%  "

$ws.Range("E43").Value = "  +9.08%  "

$ws.Range("E44").Value = "  +3.55%  "

$ws.Range("D45").Value = "2.965.84"
$ws.Range("E45").Value = "  +1.60%  "

$ws.Range("E46").Value = "  +3.93%  "

$ws.Range("D47").Value = "'38.67"
$ws.Range("E47").Value = "  +5.15%  "

$ws.Range("D48").Value = "'131.94"
$ws.Range("E48").Value = "  +4.89%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'25.80"
$ws.Range("E49").Value = "  +5.41%  "

$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "'2.34"
$ws.Range("E50").Value = "  +7.46%  "

$ws.Range("E51").Value = "  +0.01%  "
